$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) — column F "想去人数" (want-to-go count) bumps
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1108
$ws1.Range("F3").Value = 4789
$ws1.Range("F5").Value = 1935
$ws1.Range("F6").Value = 576
$ws1.Range("F12").Value = 851
$ws1.Range("F14").Value = 2026
$ws1.Range("F15").Value = 637
$ws1.Range("F20").Value = 125
$ws1.Range("F21").Value = 125
$ws1.Range("F34").Value = 4410

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F26").Value = 210
$ws2.Range("F35").Value = 50

# Sheet "本地生活" (local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F7").Value = 433

# Sheet "全部类型" (all types) — combined roll-up sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 433
$ws4.Range("F7").Value = 1108
$ws4.Range("F8").Value = 4789
$ws4.Range("F9").Value = 1935
$ws4.Range("F10").Value = 576
$ws4.Range("F20").Value = 851
$ws4.Range("F22").Value = 2026
$ws4.Range("F23").Value = 637
$ws4.Range("F29").Value = 125
$ws4.Range("F30").Value = 125
$ws4.Range("F49").Value = 4410
$ws4.Range("F50").Value = 50
